{"js": "// Replace the \"registrare le informazioni di un libro/giornale/multimedia...\"\n// list item's text with the new \"registrare gli eventi...\" text, and insert a\n// brand-new list item (same list/paragraph formatting) right after it that\n// keeps the original text.\n\nconst OLD_TEXT =\n  \"Possibilit\u00e0 di registrare le informazioni di un libro/giornale/multimedia non ancora registrato da nessuno e di registrare le informazioni della copia del libro/giornale/multimedia di cui la biblioteca \u00e8 in possesso.\";\nconst NEW_TEXT =\n  \"Possibilit\u00e0 di registrare gli eventi organizzati dalla biblioteca presso cui il bibliotecario lavora.\";\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === OLD_TEXT) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the target paragraph with the original text.\");\n}\n\n// Insert a new paragraph right after the target; it inherits the target's\n// paragraph formatting (list style/numbering), and seed it with the original\n// text.\ntarget.insertParagraph(OLD_TEXT, \"After\");\n\n// Now overwrite the original paragraph's own text with the new wording.\ntarget.insertText(NEW_TEXT, \"Replace\");\n\nawait context.sync();\n", "ps1": "# Replace the \"registrare le informazioni di un libro/giornale/multimedia...\"\n# list item's text with the new \"registrare gli eventi...\" text, and insert a\n# brand-new list item (same list/paragraph formatting) right after it that\n# keeps the original text.\n\n$d = $word.ActiveDocument\n\n$OLD_TEXT = \"Possibilit\u00e0 di registrare le informazioni di un libro/giornale/multimedia non ancora registrato da nessuno e di registrare le informazioni della copia del libro/giornale/multimedia di cui la biblioteca \u00e8 in possesso.\"\n$NEW_TEXT = \"Possibilit\u00e0 di registrare gli eventi organizzati dalla biblioteca presso cui il bibliotecario lavora.\"\n\n$target = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    # Paragraph.Range.Text includes the trailing paragraph mark (chr 13);\n    # strip it before comparing against the plain sentence text.\n    $t = $p.Range.Text.TrimEnd([char]13)\n    if ($t -eq $OLD_TEXT) {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find the target paragraph with the original text.\"\n}\n\n# Split a new (empty) paragraph in right after the target; it inherits the\n# target's paragraph formatting (list style/numbering: Paragrafoelenco /\n# numId 1) automatically, same as Word does when you press Enter.\n$target.Range.InsertParagraphAfter()\n$newPara = $target.Next()\n\n# Seed the new paragraph with the original wording (InsertBefore keeps the\n# paragraph's own mark intact, unlike assigning straight to .Range.Text on a\n# mark-only range).\n$newPara.Range.InsertBefore($OLD_TEXT)\n\n# Finally, overwrite the original paragraph's own text with the new wording.\n$target.Range.Text = $NEW_TEXT\n"}
